$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.004833799198163724
$ws.Range("E2").Value = 0.004833799198163724

# Row 3
$ws.Range("D3").Value = 0.009997074894470999
$ws.Range("E3").Value = 0.009997074894470999

# Row 4
$ws.Range("D4").Value = [double]"1.086933051355286E-07"
$ws.Range("E4").Value = [double]"1.086933051355286E-07"

# Row 5
$ws.Range("D5").Value = 0.00156078344145532
$ws.Range("E5").Value = 0.00156078344145532

# Row 6
$ws.Range("D6").Value = 0.1025729809860371
$ws.Range("E6").Value = 0.1025729809860371

# Row 7
$ws.Range("D7").Value = 0.8621228510711392
$ws.Range("E7").Value = 0.1378771489288608

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"2.539386001995285E-07"
$ws.Range("E8").Value = 0.9999997460613999

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = [double]"5.650572907967871E-08"
$ws.Range("E9").Value = 0.9999999434942709

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"4.345717238303671E-05"
$ws.Range("E10").Value = 0.9999565428276169

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.1409516226474173
$ws.Range("E11").Value = 0.8590483773525828
$ws.Range("F11").Value = 4.415120601654053
$ws.Range("G11").Value = 0.6

# Row 12
$ws.Range("D12").Value = 0.3221736294205497
$ws.Range("E12").Value = 0.3221736294205497

# Row 13
$ws.Range("D13").Value = 0.04121893225396613
$ws.Range("E13").Value = 0.04121893225396613

# Row 14
$ws.Range("D14").Value = [double]"2.598808055917455E-11"
$ws.Range("E14").Value = [double]"2.598808055917455E-11"

# Row 15
$ws.Range("D15").Value = 0.0001361609622577465
$ws.Range("E15").Value = 0.0001361609622577465

# Row 16
$ws.Range("D16").Value = 0.03786170213428709
$ws.Range("E16").Value = 0.03786170213428709

# Row 17
$ws.Range("D17").Value = 0.9412671788327128
$ws.Range("E17").Value = 0.05873282116728717

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"4.870346739783319E-11"
$ws.Range("E18").Value = 0.9999999999512965

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"1.069533225220047E-07"
$ws.Range("E19").Value = 0.9999998930466775

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"1.23269073670916E-06"
$ws.Range("E20").Value = 0.9999987673092633

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.04500703619786998
$ws.Range("E21").Value = 0.9549929638021301
$ws.Range("F21").Value = 5.703360557556152
$ws.Range("G21").Value = 0.6

$wb.Save()
